$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.6666666666666666
$ws.Range("F2").Value = 0.8333333333333334
$ws.Range("G2").Value = 0.9629629629629629
$ws.Range("H2").Value = 0.7761278037284857
$ws.Range("I2").Value = 534
$ws.Range("J2").Value = 534
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# ---- Sheet: Classification Report ----
$ws = $wb.Worksheets.Item("Classification Report")
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0.5
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.6666666666666666

$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.5

$ws.Range("B5").Value = 0.25
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 0.3333333333333333

$ws.Range("B6").Value = 0.25
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 0.3333333333333333

# ---- Sheet: Confusion Matrix ----
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 534

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 534
